$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "58.075.65"
$ws.Range("E2").Value = "  +2.15%  "

$ws.Range("D3").Value = "2.355.49"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("E4").Value = "  -0.18%  "

Set-TextValue "D5" "540.99"
$ws.Range("E5").Value = "  +2.33%  "

Set-TextValue "D6" "136.07"
$ws.Range("E6").Value = "  +2.89%  "

$ws.Range("E7").Value = "  +0.32%  "

Set-TextValue "D8" "0.565"
$ws.Range("E8").Value = "  +5.73%  "

Set-TextValue "D9" "0.103"
$ws.Range("E9").Value = "  +1.49%  "

$ws.Range("E10").Value = "  +4.78%  "

$ws.Range("E11").Value = "  -0.71%  "

$ws.Range("E12").Value = "  +2.16%  "

Set-TextValue "D13" "23.89"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("D14").Value = "2.773.16"
$ws.Range("E14").Value = "  +1.48%  "

$ws.Range("D15").Value = "58.021.23"
$ws.Range("E15").Value = "  +2.00%  "

$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("D17").Value = "2.339.40"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("E18").Value = "  +3.01%  "

Set-TextValue "D19" "332.49"
$ws.Range("E19").Value = "  -0.79%  "

Set-TextValue "D20" "4.28"
$ws.Range("E20").Value = "  +2.91%  "

Set-TextValue "D21" "6.74"
$ws.Range("E21").Value = "  -1.06%  "

Set-TextValue "D22" "1.00"
$ws.Range("E22").Value = "  +0.32%  "

Set-TextValue "D23" "62.97"
$ws.Range("E23").Value = "  +1.94%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D26" "8.50"
$ws.Range("E26").Value = "  -1.96%  "

$ws.Range("E27").Value = "  +1.61%  "

Set-TextValue "D28" "172.45"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("E29").Value = "  +2.12%  "

$ws.Range("D30").Value = "0.0₃0738"
$ws.Range("E30").Value = "  +2.07%  "

$ws.Range("E31").Value = "  +0.97%  "

Set-TextValue "D32" "1.04"
$ws.Range("E32").Value = "  +11.84%  "

Set-TextValue "D33" "18.50"
$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("E36").Value = "  +6.16%  "

Set-TextValue "D37" "1.26"
$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("E38").Value = "  +4.32%  "

Set-TextValue "D39" "39.30"
$ws.Range("E39").Value = "  +0.50%  "

Set-TextValue "D40" "145.32"
$ws.Range("E40").Value = "  -1.99%  "

Set-TextValue "D41" "295.29"
$ws.Range("E41").Value = "  +4.57%  "

$ws.Range("E42").Value = "  +0.91%  "

Set-TextValue "D43" "3.65"
$ws.Range("E43").Value = "  +1.50%  "

Set-TextValue "D44" "0.0949"
$ws.Range("E44").Value = "  +1.89%  "

Set-TextValue "D45" "19.34"
$ws.Range("E45").Value = "  +3.03%  "

Set-TextValue "D46" "0.0503"
$ws.Range("E46").Value = "  +0.84%  "

Set-TextValue "D47" "0.565"
$ws.Range("E47").Value = "  +1.22%  "

$ws.Range("E48").Value = "  +2.50%  "

Set-TextValue "D49" "0.382"
$ws.Range("E49").Value = "  -0.03%  "

Set-TextValue "D50" "17.47"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("E51").Value = "  +0.48%  "
